# Generate Report for Handback
# Row 7 (80b56b9a-fbb9-44d5-8df7-0f1b524176c2.md) now has a handback generated for it,
# but the handback version doesn't match the latest handoff, so an error is recorded.

$wb = $excel.ActiveWorkbook

$targetUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/7411f0d8f238462910f04ad9fc202205c5b37ba1/e2e/80b56b9a-fbb9-44d5-8df7-0f1b524176c2.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/4829b6f1601b0ab426d7832e89e79857d19be4ed/e2e/80b56b9a-fbb9-44d5-8df7-0f1b524176c2.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/7411f0d8f238462910f04ad9fc202205c5b37ba1/e2e/80b56b9a-fbb9-44d5-8df7-0f1b524176c2.md."

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("I7").Value = "80b56b9a-fbb9-44d5-8df7-0f1b524176c2.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I7"), $targetUrl, [Type]::Missing, [Type]::Missing, "80b56b9a-fbb9-44d5-8df7-0f1b524176c2.md")
$wsZhCn.Range("I7").Style = "HyperLink"

$wsZhCn.Range("J7").Value = "80b56b9a-fbb9-44d5-8df7-0f1b524176c2.3f38b1cdf108a94d7d16a47dfeee58922cb8e645.zh-cn.xlf"
$wsZhCn.Range("K7").Value = "2016-08-13 01:07:59"
$wsZhCn.Range("P7").Value = $errorDetail

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("I7").Value = "80b56b9a-fbb9-44d5-8df7-0f1b524176c2.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I7"), $targetUrl, [Type]::Missing, [Type]::Missing, "80b56b9a-fbb9-44d5-8df7-0f1b524176c2.md")
$wsDeDe.Range("I7").Style = "HyperLink"

$wsDeDe.Range("J7").Value = "80b56b9a-fbb9-44d5-8df7-0f1b524176c2.3f38b1cdf108a94d7d16a47dfeee58922cb8e645.de-de.xlf"
$wsDeDe.Range("K7").Value = "2016-08-13 01:08:13"
$wsDeDe.Range("P7").Value = $errorDetail
